# Insert a new data row at row 106 (pushing the existing rows 106-137
# down to 107-138) and populate it with the new weekly price entry for
# Membrillo at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 106..137 down to 107..138, leaving a blank row 106.
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the inserted record.
$ws.Cells.Item(106, 1).Value  = 4
$ws.Cells.Item(106, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(106, 3).Value  = "Los Lagos"
$ws.Cells.Item(106, 4).Value  = 45029
$ws.Cells.Item(106, 5).Value  = 10
$ws.Cells.Item(106, 6).Value  = "Fruta"
$ws.Cells.Item(106, 7).Value  = 100104
$ws.Cells.Item(106, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(106, 9).Value  = 100104003
$ws.Cells.Item(106, 10).Value = "Membrillo"
$ws.Cells.Item(106, 11).Value = "Champion"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 200
$ws.Cells.Item(106, 14).Value = 15000
$ws.Cells.Item(106, 15).Value = 16000
$ws.Cells.Item(106, 16).Value = 15500
$ws.Cells.Item(106, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(106, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(106, 19).Value = 861
$ws.Cells.Item(106, 20).Value = 18
